$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# ColumnWidth (character units) maps to the OOXML "width" attribute via Excel's own
# pixel-quantization rounding. Empirically, an input of (target - 0.8335) reliably
# lands in the middle of the input range that stores as the integer target width.
$ws.Columns.Item(13).ColumnWidth = 19.1665   # col 13 (M): 9 -> 20
$ws.Columns.Item(19).ColumnWidth = 18.1665   # col 19 (S): 20 -> 19
$ws.Columns.Item(20).ColumnWidth = 19.1665   # col 20 (T): 15 -> 20
$ws.Columns.Item(24).ColumnWidth = 18.1665   # col 24 (X): 13 -> 19

# --- Cell B2 ---
$ws.Range("B2").Value = "work"

# --- Row 3 "Total Cost ($)" unit-price values (Len+Units incorporated) ---
$ws.Range("L3").Value = 181.44
$ws.Range("M3").Value = 6.48
$ws.Range("N3").Value = 19.04
$ws.Range("O3").Value = 10.248
$ws.Range("P3").Value = 5.598
$ws.Range("Q3").Value = 1.944
$ws.Range("R3").Value = 36
$ws.Range("S3").Value = 9.234000000000002
$ws.Range("T3").Value = 12.312
$ws.Range("U3").Value = 11.76
$ws.Range("V3").Value = 6.804
$ws.Range("W3").Value = 15.24
$ws.Range("X3").Value = 279.3333333333333
$ws.Range("Y3").Value = 209.5
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 497
$ws.Range("AB3").Value = 263.5
$ws.Range("AC3").Value = 222
$ws.Range("AD3").Value = 132
$ws.Range("AE3").Value = 0
